$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" (before "2021-Q4").
#    Copying the existing "2021-Q4" sheet gives us an identical template
#    (sheetPr / header row style / page margins) to build the new data on.
# ---------------------------------------------------------------------------
$total    = $wb.Worksheets.Item(1)
$template = $wb.Worksheets.Item(2)
$template.Copy($null, $total)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the fund holdings table.
#    (Header row B1:H1 is already correct, copied from the template.)
# ---------------------------------------------------------------------------
$data = @(
    @("013385", "信澳优势价值混合A", "12.44", "84.28", "3.70", "0.4603", 9),
    @("013393", "信澳价值精选混合A", "3.34", "79.98", "3.51", "0.1172", 7),
    @("003655", "信澳新财富灵活配置混合", "4.04", "54.87", "1.73", "0.0699", 8),
    @("013386", "信澳优势价值混合C", "1.26", "84.28", "3.70", "0.0466", 9),
    @("673090", "西部利得个股精选股票A", "1.25", "86.69", "2.47", "0.0309", 7),
    @("013554", "信澳远见价值混合A", "0.92", "48.39", "2.89", "0.0266", 5),
    @("013262", "西部利得个股精选股票C", "0.95", "86.69", "2.47", "0.0235", 7),
    @("013555", "信澳远见价值混合C", "0.64", "48.39", "2.89", "0.0185", 5),
    @("013394", "信澳价值精选混合C", "0.38", "79.98", "3.51", "0.0133", 7),
    @("162211", "泰达宏利品质生活混合", "0.11", "71.38", "3.13", "0.0034", 10),
    @("015694", "瑞达策略优选混合A", "0.09", "67.87", "2.76", "0.0025", 7),
    @("015695", "瑞达策略优选混合C", "0.00", "67.87", "2.76", 0, 7)
)

# The template copy only had 2 data rows (row 2/3), both already carrying
# the correct column-A style (s="2"). Stamp that same style across the rest
# of column A (rows 4-13) before writing the real index values into it.
$newWs.Range("A2").Copy($newWs.Range("A4:A13"))

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    # Column A (index) is a genuine number and needs no special formatting.
    $newWs.Cells.Item($r, 1).Value = $i

    # Columns B-F are text-looking numbers stored as real text (inlineStr in
    # the source file) - force text format, write, then drop the format
    # residue so the cell keeps the default (unstyled) look.
    $newWs.Cells.Item($r, 2).NumberFormat = "@"
    $newWs.Cells.Item($r, 2).Value = $row[0]
    $newWs.Cells.Item($r, 2).ClearFormats()

    $newWs.Cells.Item($r, 3).NumberFormat = "@"
    $newWs.Cells.Item($r, 3).Value = $row[1]
    $newWs.Cells.Item($r, 3).ClearFormats()

    $newWs.Cells.Item($r, 4).NumberFormat = "@"
    $newWs.Cells.Item($r, 4).Value = $row[2]
    $newWs.Cells.Item($r, 4).ClearFormats()

    $newWs.Cells.Item($r, 5).NumberFormat = "@"
    $newWs.Cells.Item($r, 5).Value = $row[3]
    $newWs.Cells.Item($r, 5).ClearFormats()

    $newWs.Cells.Item($r, 6).NumberFormat = "@"
    $newWs.Cells.Item($r, 6).Value = $row[4]
    $newWs.Cells.Item($r, 6).ClearFormats()

    if ($row[5] -eq 0) {
        # Row 13's "持有市值" is a genuine number 0, not a text cell.
        $newWs.Cells.Item($r, 7).Value = 0
    } else {
        $newWs.Cells.Item($r, 7).NumberFormat = "@"
        $newWs.Cells.Item($r, 7).Value = $row[5]
        $newWs.Cells.Item($r, 7).ClearFormats()
    }

    # Column H (rank) is a genuine number and needs no special formatting.
    $newWs.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new row right under the
#    header for the 2022-Q3 figures and shift the existing rows down.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").Clear()

# Copy column-A style (s="2") from the row below onto the new row, then set
# the real values.
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 0.81

# Renumber the index column (A) for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 4. Restore the originally active sheet/tab (last sheet, "2021-Q1").
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()

Write-Output "edit complete"
